$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B (Coin) and C (Link) are plain text; D (Price) and E (Volume) must
# be forced to text so numeric-looking strings are not converted to numbers.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 2
Set-TextValue $ws.Range('D2') '42.180.87'
Set-TextValue $ws.Range('E2') '  -1.37%  '

# Row 3
Set-TextValue $ws.Range('D3') '2.502.26'
Set-TextValue $ws.Range('E3') '  -2.67%  '

# Row 4
Set-TextValue $ws.Range('E4') '  +0.07%  '

# Row 5
Set-TextValue $ws.Range('D5') '302.20'
Set-TextValue $ws.Range('E5') '  +0.41%  '

# Row 6
Set-TextValue $ws.Range('D6') '95.54'
Set-TextValue $ws.Range('E6') '  -1.24%  '

# Row 7
Set-TextValue $ws.Range('D7') '0.582'
Set-TextValue $ws.Range('E7') '  +1.52%  '

# Row 8
Set-TextValue $ws.Range('E8') '  +0.04%  '

# Row 9
Set-TextValue $ws.Range('D9') '0.533'
Set-TextValue $ws.Range('E9') '  -2.58%  '

# Row 10
Set-TextValue $ws.Range('D10') '35.93'
Set-TextValue $ws.Range('E10') '  -0.80%  '

# Row 11
Set-TextValue $ws.Range('D11') '0.0807'
Set-TextValue $ws.Range('E11') '  -0.06%  '

# Row 12
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Range('D12') '7.59'
Set-TextValue $ws.Range('E12') '  -1.00%  '

# Row 13
$ws.Range('B13').Value = 'TRON'
$ws.Range('C13').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-TextValue $ws.Range('D13') '0.112'
Set-TextValue $ws.Range('E13') '  -2.81%  '

# Row 14
Set-TextValue $ws.Range('D14') '2.883.92'
Set-TextValue $ws.Range('E14') '  -2.79%  '

# Row 15
Set-TextValue $ws.Range('D15') '2.542.05'
Set-TextValue $ws.Range('E15') '  +3.49%  '

# Row 16
Set-TextValue $ws.Range('D16') '14.99'
Set-TextValue $ws.Range('E16') '  +4.26%  '

# Row 17
Set-TextValue $ws.Range('D17') '0.850'
Set-TextValue $ws.Range('E17') '  -3.38%  '

# Row 18
Set-TextValue $ws.Range('D18') '42.166.99'
Set-TextValue $ws.Range('E18') '  -1.48%  '

# Row 19
Set-TextValue $ws.Range('D19') '12.71'
Set-TextValue $ws.Range('E19') '  -0.91%  '

# Row 20
Set-TextValue $ws.Range('E20') '  -2.29%  '

# Row 21
Set-TextValue $ws.Range('D21') '6.36'
Set-TextValue $ws.Range('E21') '  -3.84%  '

# Row 22
Set-TextValue $ws.Range('D22') '70.82'
Set-TextValue $ws.Range('E22') '  -1.58%  '

# Row 23
Set-TextValue $ws.Range('D23') '248.47'
Set-TextValue $ws.Range('E23') '  -2.24%  '

# Row 24
Set-TextValue $ws.Range('E24') '  -1.85%  '

# Row 25
Set-TextValue $ws.Range('D25') '2.00'
Set-TextValue $ws.Range('E25') '  -4.85%  '

# Row 26
Set-TextValue $ws.Range('D26') '26.84'
Set-TextValue $ws.Range('E26') '  -5.11%  '

# Row 27
Set-TextValue $ws.Range('D27') '1.00'
Set-TextValue $ws.Range('E27') '  +0.09%  '

# Row 28
Set-TextValue $ws.Range('D28') '2.31'
Set-TextValue $ws.Range('E28') '  +9.75%  '

# Row 29
Set-TextValue $ws.Range('D29') '10.07'
Set-TextValue $ws.Range('E29') '  -1.05%  '

# Row 30
Set-TextValue $ws.Range('D30') '36.95'
Set-TextValue $ws.Range('E30') '  -5.34%  '

# Row 31
Set-TextValue $ws.Range('D31') '5.86'
Set-TextValue $ws.Range('E31') '  -1.67%  '

# Row 32
Set-TextValue $ws.Range('D32') '153.34'
Set-TextValue $ws.Range('E32') '  -1.33%  '

# Row 33
Set-TextValue $ws.Range('D33') '3.24'
Set-TextValue $ws.Range('E33') '  -3.03%  '

# Row 34
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue $ws.Range('D34') '0.0778'
Set-TextValue $ws.Range('E34') '  -3.65%  '

# Row 35
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range('D35') '2.05'
Set-TextValue $ws.Range('E35') '  -4.93%  '

# Row 36
$ws.Range('B36').Value = 'Celestia'
$ws.Range('C36').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextValue $ws.Range('D36') '18.49'
Set-TextValue $ws.Range('E36') '  -0.68%  '

# Row 37
$ws.Range('B37').Value = 'WEMIXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-TextValue $ws.Range('D37') '2.59'
Set-TextValue $ws.Range('E37') '  -5.72%  '

# Row 38
Set-TextValue $ws.Range('D38') '0.114'
Set-TextValue $ws.Range('E38') '  +0.21%  '

# Row 39
Set-TextValue $ws.Range('D39') '0.119'
Set-TextValue $ws.Range('E39') '  +0.11%  '

# Row 40
Set-TextValue $ws.Range('D40') '23.53'
Set-TextValue $ws.Range('E40') '  -0.03%  '

# Row 41
$ws.Range('B41').Value = 'NEARProtocol'
$ws.Range('C41').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue $ws.Range('D41') '3.35'
Set-TextValue $ws.Range('E41') '  -1.26%  '

# Row 42
$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range('D42') '3.80'
Set-TextValue $ws.Range('E42') '  -1.27%  '

# Row 43
Set-TextValue $ws.Range('D43') '1.00'
Set-TextValue $ws.Range('E43') '  +0.02%  '

# Row 44
Set-TextValue $ws.Range('D44') '2.017.99'
Set-TextValue $ws.Range('E44') '  -2.79%  '

# Row 45
Set-TextValue $ws.Range('D45') '0.0295'
Set-TextValue $ws.Range('E45') '  -4.01%  '

# Row 46
Set-TextValue $ws.Range('D46') '1.92'
Set-TextValue $ws.Range('E46') '  -9.93%  '

# Row 47
Set-TextValue $ws.Range('D47') '83.49'
Set-TextValue $ws.Range('E47') '  -1.04%  '

# Row 48
Set-TextValue $ws.Range('D48') '8.93'
Set-TextValue $ws.Range('E48') '  -3.90%  '

# Row 49
Set-TextValue $ws.Range('D49') '2.743.18'
Set-TextValue $ws.Range('E49') '  -2.79%  '

# Row 50
Set-TextValue $ws.Range('D50') '71.92'
Set-TextValue $ws.Range('E50') '  -5.87%  '

# Row 51
$ws.Range('B51').Value = 'Aave'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range('D51') '100.43'
Set-TextValue $ws.Range('E51') '  -4.55%  '
